# Electronics/MCU_details_and_scoring.xlsx
# "changed pro mini link and price"
#
# The Arduino PRO Mini row (row 4) had its per-piece price recomputed as a
# text formula ("$" & ROUND(9.99/3,2) & " per piece" => "$3.33 per piece").
# Replace it with the new, plain numeric unit price (12.99) formatted as
# Currency, and leave the cursor where the author left it (G13).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$priceCell = $ws.Range("E4")

# Drop the old formula and store the new price as a literal currency value.
$priceCell.Value = 12.99
$priceCell.NumberFormat = "_(""$""* #,##0.00_);_(""$""* \(#,##0.00\);_(""$""* ""-""??_);_(@_)"

# Restore the author's last selection.
$ws.Range("G13").Select() | Out-Null
